$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text so values like "67.404.77" or "597.47"
# are not auto-converted to dates/numbers by Excel type inference.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '67.404.77'
$ws.Range('E2').Value = '  +1.15%  '

$ws.Range('D3').Value = '3.525.92'
$ws.Range('E3').Value = '  +0.66%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '597.47'
$ws.Range('E5').Value = '  +1.23%  '

$ws.Range('D6').Value = '173.50'
$ws.Range('E6').Value = '  +2.53%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  +1.66%  '

$ws.Range('E9').Value = '  +8.56%  '

$ws.Range('D10').Value = '7.31'
$ws.Range('E10').Value = '  +0.61%  '

$ws.Range('D11').Value = '0.436'
$ws.Range('E11').Value = '  -0.01%  '

$ws.Range('D12').Value = '4.134.07'
$ws.Range('E12').Value = '  +0.59%  '

$ws.Range('E13').Value = '  -0.17%  '

$ws.Range('D14').Value = '28.75'
$ws.Range('E14').Value = '  +2.80%  '

$ws.Range('D15').Value = '0.0000182'
$ws.Range('E15').Value = '  +2.28%  '

$ws.Range('D16').Value = '67.297.31'
$ws.Range('E16').Value = '  +1.01%  '

$ws.Range('D17').Value = '3.519.56'
$ws.Range('E17').Value = '  +0.60%  '

$ws.Range('D18').Value = '6.36'
$ws.Range('E18').Value = '  +1.23%  '

$ws.Range('D19').Value = '14.20'
$ws.Range('E19').Value = '  +1.03%  '

$ws.Range('D20').Value = '397.09'
$ws.Range('E20').Value = '  +2.19%  '

$ws.Range('D21').Value = '7.99'
$ws.Range('E21').Value = '  +0.18%  '

$ws.Range('D22').Value = '73.71'
$ws.Range('E22').Value = '  +1.00%  '

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '0.541'
$ws.Range('E23').Value = '  +2.42%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('D25').Value = '0.0000124'
$ws.Range('E25').Value = '  +0.28%  '

$ws.Range('D26').Value = '10.26'
$ws.Range('E26').Value = '  +0.57%  '

$ws.Range('E27').Value = '  +0.48%  '

$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.09%  '

$ws.Range('D29').Value = '6.31'
$ws.Range('E29').Value = '  -0.59%  '

$ws.Range('D30').Value = '1.47'
$ws.Range('E30').Value = '  -0.07%  '

$ws.Range('D31').Value = '2.09'
$ws.Range('E31').Value = '  +1.51%  '

$ws.Range('D32').Value = '24.10'
$ws.Range('E32').Value = '  +2.71%  '

$ws.Range('D33').Value = '7.38'
$ws.Range('E33').Value = '  -0.30%  '

$ws.Range('D34').Value = '1.64'
$ws.Range('E34').Value = '  +5.34%  '

$ws.Range('D35').Value = '164.45'
$ws.Range('E35').Value = '  +2.07%  '

$ws.Range('D36').Value = '0.899'
$ws.Range('E36').Value = '  -0.65%  '

$ws.Range('D37').Value = '1.92'
$ws.Range('E37').Value = '  -0.18%  '

$ws.Range('D38').Value = '4.74'
$ws.Range('E38').Value = '  +1.62%  '

$ws.Range('D39').Value = '6.86'
$ws.Range('E39').Value = '  +2.24%  '

$ws.Range('D40').Value = '0.0748'
$ws.Range('E40').Value = '  +0.23%  '

$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.66'
$ws.Range('E41').Value = '  +4.28%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '26.54'
$ws.Range('E42').Value = '  +0.59%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '27.19'
$ws.Range('E43').Value = '  +0.75%  '

$ws.Range('D44').Value = '2.819.13'
$ws.Range('E44').Value = '  +1.09%  '

$ws.Range('D45').Value = '43.02'
$ws.Range('E45').Value = '  -1.20%  '

$ws.Range('D46').Value = '0.0312'
$ws.Range('E46').Value = '  -0.91%  '

$ws.Range('D47').Value = '341.52'
$ws.Range('E47').Value = '  -4.66%  '

$ws.Range('D48').Value = '1.10'
$ws.Range('E48').Value = '  +0.88%  '

$ws.Range('D49').Value = '33.81'
$ws.Range('E49').Value = '  +2.28%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '6.52'
$ws.Range('E50').Value = '  +0.48%  '

$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '0.854'
$ws.Range('E51').Value = '  +0.31%  '

# Restore the original (default) style so the cell styling/number
# format matches the source workbook exactly (no lingering "@" text format).
$dRange.Style = "Normal"
